# Adds the new "buyItem" event-chain rows (138-149) to eventAction.csv sheet
# and updates the active selection to reflect the author's final cursor position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 138
$ws.Range('A138').Value = 'buyItem'
$ws.Range('B138').Value = '购买道具'
$ws.Range('C138').Value = 'eventList'
$ws.Range('D138').Value = 'close;cityHaveItem'

# Row 139
$ws.Range('A139').Value = 'cityHaveItem'
$ws.Range('C139').Value = 'condition'
$ws.Range('D139').Value = 'cityHaveItem;buyItemStart;buyItemNoItem'

# Row 140
$ws.Range('A140').Value = 'buyItemStart'
$ws.Range('C140').Value = 'eventList'
$ws.Range('D140').Value = 'buyItemDialog;buyItemWindow'

# Row 141
$ws.Range('A141').Value = 'buyItemDialog'
$ws.Range('C141').Value = 'dialog'
$ws.Range('D141').Value = 'dialog_what_do_you_wanna_buy'

# Row 142
$ws.Range('A142').Value = 'buyItemNoItem'
$ws.Range('C142').Value = 'eventList'
$ws.Range('D142').Value = 'buyItemNoItemDialog;shop'

# Row 143
$ws.Range('A143').Value = 'buyItemNoItemDialog'
$ws.Range('C143').Value = 'dialog'
$ws.Range('D143').Value = 'dialog_no_item_to_sell'

# Row 144
$ws.Range('A144').Value = 'buyItemWindow'
$ws.Range('C144').Value = 'window'
$ws.Range('D144').Value = 'ItemBrowsePanel;buyItemSelected;buyItemCanceled;0'

# Row 145
$ws.Range('A145').Value = 'buyItemSelected'
$ws.Range('C145').Value = 'eventList'
$ws.Range('D145').Value = 'buyItemSmallWindow'

# Row 146
$ws.Range('A146').Value = 'buyItemCanceled'
$ws.Range('C146').Value = 'eventList'
$ws.Range('D146').Value = 'shop'

# Row 147
$ws.Range('A147').Value = 'buyItemSmallWindow'
$ws.Range('C147').Value = 'window'
$ws.Range('D147').Value = 'ItemInfoPanel;buyItemConfirmed;buyItemCancelSmallWindow;0;reserved.itemId'

# Row 148
$ws.Range('A148').Value = 'buyItemConfirmed'
$ws.Range('C148').Value = 'eventList'
$ws.Range('D148').Value = ';'

# Row 149
$ws.Range('A149').Value = 'buyItemCancelSmallWindow'
$ws.Range('C149').Value = 'eventList'
$ws.Range('D149').Value = ';'

# Column C cells use the eventAction "keyword" style (Chinese-capable font, style index 1).
# Copy number/font formatting only (PasteSpecial xlPasteFormats) from an existing cell that
# already carries the matching keyword text, so the style is reused instead of duplicated.
$ws.Range('C19').Copy() | Out-Null
$ws.Range('C138').PasteSpecial(-4122) | Out-Null
$ws.Range('C36').Copy() | Out-Null
$ws.Range('C139').PasteSpecial(-4122) | Out-Null
$ws.Range('C19').Copy() | Out-Null
$ws.Range('C140').PasteSpecial(-4122) | Out-Null
$ws.Range('C16').Copy() | Out-Null
$ws.Range('C141').PasteSpecial(-4122) | Out-Null
$ws.Range('C19').Copy() | Out-Null
$ws.Range('C142').PasteSpecial(-4122) | Out-Null
$ws.Range('C16').Copy() | Out-Null
$ws.Range('C143').PasteSpecial(-4122) | Out-Null
$ws.Range('C22').Copy() | Out-Null
$ws.Range('C144').PasteSpecial(-4122) | Out-Null
$ws.Range('C19').Copy() | Out-Null
$ws.Range('C145').PasteSpecial(-4122) | Out-Null
$ws.Range('C19').Copy() | Out-Null
$ws.Range('C146').PasteSpecial(-4122) | Out-Null
$ws.Range('C22').Copy() | Out-Null
$ws.Range('C147').PasteSpecial(-4122) | Out-Null
$ws.Range('C19').Copy() | Out-Null
$ws.Range('C148').PasteSpecial(-4122) | Out-Null
$ws.Range('C19').Copy() | Out-Null
$ws.Range('C149').PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the saved selection/scroll position to match the author's final view (D144 active).
$ws.Activate()
$ws.Range('D144').Select()
